{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"Drawers opening: / closing\" paragraph so we can insert the\n// new \"Cupboard Door Open / Close\" item right after it (and before the\n// \"Wardrobe\" item).\nlet drawersPara = null;\nlet walkingPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"Drawers opening\") !== -1) {\n    drawersPara = paragraphs.items[i];\n  }\n  if (t.indexOf(\"Walking sounds\") !== -1) {\n    walkingPara = paragraphs.items[i];\n  }\n}\n\nif (!drawersPara || !walkingPara) {\n  throw new Error(\"Could not locate anchor paragraphs\");\n}\n\n// Insert the new \"Cupboard Door Open / Close\" list item right after\n// \"Drawers opening: / closing\". It inherits the ListParagraph style /\n// numbering from the paragraph it is inserted relative to.\nconst cupboardPara = drawersPara.insertParagraph(\n  \"Cupboard Door Open / Close\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// The document's \"_GoBack\" bookmark used to sit at the end of the\n// \"Walking sounds: wood walk\" paragraph; it now belongs in the middle of\n// the new paragraph, between \"Cupboard Door Open / \" and \"Close\". Remove\n// the old one first, then re-insert it in its new location.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst closeResults = cupboardPara.search(\"Close\", { matchCase: true });\ncloseResults.load(\"items\");\nawait context.sync();\n\nconst closeStart = closeResults.items[0].getRange(\"Start\");\ncloseStart.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// Append the rest of the new sound-list entries after the\n// \"Walking sounds: wood walk\" paragraph, before the trailing blank\n// paragraph / section break.\nconst newItems = [\n  \"Objects: Hitting Surfaces, breaking porcelain\",\n  \"Tap on/of\",\n  \"Shower On/Off/Shower Running Sound\",\n  \"Clothes Rustling\",\n  \"Fridge Open/Close/Running (Humming Sound)\",\n  \"Toaster Push Down/Pop Up\",\n  \"Kettle On/Off/Running\",\n  \"Stove On/Off/Gas Running\",\n  \"Oven Running\",\n  \"Stirring Tea\",\n  \"Spoon Hitting Surface/Tea Mug\",\n  \"Shower Door Open/Close\",\n  \"Items Brushing across Kitchen Surface\"\n];\n\nlet anchor = walkingPara;\nfor (const itemText of newItems) {\n  anchor = anchor.insertParagraph(itemText, Word.InsertLocation.after);\n}\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# $word.ActiveDocument is available as the document to edit.\n\n$d = $word.ActiveDocument\n\n# --- Locate the anchor paragraphs by their text -----------------------\n$drawersIndex = 0\n$walkingIndex = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*Drawers opening*\") {\n        $drawersIndex = $i\n    }\n    if ($t -like \"*Walking sounds*\") {\n        $walkingIndex = $i\n    }\n}\n\n# --- Insert \"Cupboard Door Open / Close\" right after \"Drawers opening\" -\n$drawersPara = $d.Paragraphs.Item($drawersIndex)\n$drawersPara.Range.InsertParagraphAfter() | Out-Null\n$cupboardIndex = $drawersIndex + 1\n$cupboardPara = $d.Paragraphs.Item($cupboardIndex)\n$cupboardPara.Range.Text = \"Cupboard Door Open / Close\"\n\n# --- Move the \"_GoBack\" bookmark from the end of the \"Walking sounds\" --\n# --- paragraph into the middle of the new \"Cupboard Door\" paragraph ----\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$findRange = $cupboardPara.Range.Duplicate\n$findRange.Find.Execute(\"Close\") | Out-Null\n$bmRange = $d.Range($findRange.Start, $findRange.Start)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange) | Out-Null\n\n# --- Append the remaining new sound-list entries after \"Walking sounds\" -\n# The \"Walking sounds\" paragraph shifted down by one because of the new\n# \"Cupboard Door\" paragraph inserted above it.\n$walkingIndex = $walkingIndex + 1\n\n$newItems = @(\n    \"Objects: Hitting Surfaces, breaking porcelain\",\n    \"Tap on/of\",\n    \"Shower On/Off/Shower Running Sound\",\n    \"Clothes Rustling\",\n    \"Fridge Open/Close/Running (Humming Sound)\",\n    \"Toaster Push Down/Pop Up\",\n    \"Kettle On/Off/Running\",\n    \"Stove On/Off/Gas Running\",\n    \"Oven Running\",\n    \"Stirring Tea\",\n    \"Spoon Hitting Surface/Tea Mug\",\n    \"Shower Door Open/Close\",\n    \"Items Brushing across Kitchen Surface\"\n)\n\n$idx = $walkingIndex\nforeach ($t in $newItems) {\n    $p = $d.Paragraphs.Item($idx)\n    $p.Range.InsertParagraphAfter() | Out-Null\n    $idx = $idx + 1\n    $newP = $d.Paragraphs.Item($idx)\n    $newP.Range.Text = $t\n}\n"}
